# Daily update at 8 AM UTC
# Appends a new day's data row (row 96) to the Wins Over Time sheet,
# and shifts the "final row" date-only formatting from row 95 to row 96.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 95 was the previous last row and had the special date-only
# number format; now that it's no longer last, give it the standard
# date+time format used by all other data rows.
$ws.Range("A95").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 96 becomes the new last row - populate its values and give
# it the special date-only number format previously used by row 95.
$ws.Range("A96").Value = 45835
$ws.Range("B96").Value = 407
$ws.Range("C96").Value = 403
$ws.Range("D96").Value = 418
$ws.Range("A96").NumberFormat = "YYYY-MM-DD"
